$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Select the cell at the top of the column that's about to be removed (M1),
# matching the active-cell/selection state left behind after the edit, then
# delete the entire column M. This shifts the old column N (and everything
# to its right) one column to the left, so the old N column becomes the new M.
$ws.Range("M1").Select()
$ws.Columns.Item(13).Delete()
